# Sample_File_Name_Key.xlsx edit
#
# The cellenONE_Image_File_Names column (G) for the chip-run blocks at the
# bottom of the sheet (rows 216-413) previously held long, verbose,
# per-image / per-run identifiers (often one unique string per row in a
# contiguous block). This edit collapses each block down to a single,
# short "Chip_*" / "LibraryChip_*" label per chip/run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chip 227J1 (rows 216-263)
$ws.Range("G216:G263").Value = "Chip_227J1"

# Chip 227J4 (rows 264-311)
$ws.Range("G264:G311").Value = "Chip_227J4"

# Chip 2308J3 (rows 312-359)
$ws.Range("G312:G359").Value = "Chip_2308J3"

# Chip 2308J4 (rows 360-407)
$ws.Range("G360:G407").Value = "Chip_2308J4"

# Library chip 2 (rows 408-410)
$ws.Range("G408:G410").Value = "LibraryChip_2"

# Library chip 1 (rows 411-413)
$ws.Range("G411:G413").Value = "LibraryChip_1"

# Restore the last active cell / selection used when the file was saved.
$ws.Range("J413").Select()
